$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '96.411.25'
$ws.Range('E2').Value = '  +4.48%  '
$ws.Range('D3').Value = '3.125.56'
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.04'
$ws.Range('E5').Value = '  +3.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '608.18'
$ws.Range('E6').Value = '  -0.74%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.10'
$ws.Range('E7').Value = '  +1.81%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.385'
$ws.Range('E8').Value = '  -0.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').Value = '3.114.56'
$ws.Range('E10').Value = '  +0.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.790'
$ws.Range('E11').Value = '  +1.41%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.197'
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('D13').Value = '95.928.10'
$ws.Range('E13').Value = '  +4.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000239'
$ws.Range('E14').Value = '  -1.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '33.90'
$ws.Range('E15').Value = '  +0.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.33'
$ws.Range('E16').Value = '  -1.12%  '
$ws.Range('D17').Value = '3.683.30'
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('D18').Value = '3.093.81'
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.56'
$ws.Range('E19').Value = '  -5.78%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.43'
$ws.Range('E20').Value = '  +0.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '469.14'
$ws.Range('E21').Value = '  +7.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.68'
$ws.Range('E22').Value = '  -0.86%  '
$ws.Range('B23').Value = 'PEPE'
$ws.Range('C23').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0000191'
$ws.Range('E23').Value = '  -3.20%  '
$ws.Range('B24').Value = 'Uniswap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.76'
$ws.Range('E24').Value = '  -3.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.52'
$ws.Range('E25').Value = '  -0.61%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '85.21'
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.65'
$ws.Range('E27').Value = '  +2.60%  '
$ws.Range('D28').Value = '3.264.20'
$ws.Range('E28').Value = '  +0.15%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.236'
$ws.Range('E30').Value = '  +0.69%  '
$ws.Range('B31').Value = 'Cronos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.176'
$ws.Range('E31').Value = '  -1.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.125'
$ws.Range('E32').Value = '  +1.79%  '
$ws.Range('B33').Value = 'Binance-PegBSC-USD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.998'
$ws.Range('E33').Value = '  -3.99%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.05'
$ws.Range('E34').Value = '  -1.04%  '
$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '26.15'
$ws.Range('E35').Value = '  +2.46%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.37'
$ws.Range('E36').Value = '  -7.72%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.150'
$ws.Range('E37').Value = '  -3.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '488.14'
$ws.Range('E38').Value = '  +4.89%  '
$ws.Range('B39').Value = 'PancakeSwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.87'
$ws.Range('E39').Value = '  -0.89%  '
$ws.Range('B40').Value = 'WhiteBITCoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '24.14'
$ws.Range('E40').Value = '  +1.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.437'
$ws.Range('E41').Value = '  +0.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.63'
$ws.Range('E42').Value = '  -6.66%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.23'
$ws.Range('E43').Value = '  -2.90%  '
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.17'
$ws.Range('E45').Value = '  -2.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '161.77'
$ws.Range('E46').Value = '  +1.73%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.698'
$ws.Range('E47').Value = '  +2.83%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.89'
$ws.Range('E48').Value = '  +3.40%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '43.82'
$ws.Range('E49').Value = '  +0.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.997'
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('B51').Value = 'Filecoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.33'
$ws.Range('E51').Value = '  +0.58%  '
